$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at O — shifts "Sample ID" and everything after it
# one column to the right (O->P, P->Q, ... BB->BC), matching the
# existing data/format of the column that used to be at O (inherits
# style from the column to its right, consistent with the diff where
# the new column has no special style while old header row keeps s="1").
$ws.Columns("O").Insert()

# New header for the inserted column
$ws.Cells.Item(1, 15).Value2 = "min sample length"

# Record "min sample length" = 45 for every data row (2-7)
$ws.Cells.Item(2, 15).Value2 = 45
$ws.Cells.Item(3, 15).Value2 = 45
$ws.Cells.Item(4, 15).Value2 = 45
$ws.Cells.Item(5, 15).Value2 = 45
$ws.Cells.Item(6, 15).Value2 = 45
$ws.Cells.Item(7, 15).Value2 = 45

# Rename "bucket" -> "bucket_sediment" for the bucket measurement rows
$ws.Cells.Item(2, 10).Value2 = "bucket_sediment"
$ws.Cells.Item(3, 10).Value2 = "bucket_sediment"

# Update derived Sample ID text to reflect the renamed measurement device
$ws.Cells.Item(3, 16).Value2 = "2020_07_30_17h12m20s_vault-lake_bucket_sediment_CO2"

# Row 2 (previously-rejected bucket CH4 run) now also gets its R_value_used
# (r^2) recorded, column W (23) after the insert
$ws.Cells.Item(2, 23).Value2 = 0.5393535609116455
